# Add "multiple columns primary key" support to the csv_schema sheet.
#
# This inserts a new "data_type" column (C) that holds the SQL data type of
# each column (e.g. DATETIME), separate from the "type" column which now
# only flags primary-key membership (PK). It also demonstrates a composite
# (multi-column) primary key for AutoDrive_Event_100_20201215_Sample by
# adding ActivityID and Attribute_101 as extra PK rows for that table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before C ("data_type"), shifting the old
#    type/target_table/target_column columns one position to the right.
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("C1").Value = "data_type"
$ws.Columns("C").ColumnWidth = 13.83

# 2) Every existing row that carried a data type (stored as "DATETIME" in
#    what is now column D, the old "type" column) actually belongs in the
#    new "data_type" column - move it over and clear the old spot.
$lastRow = $ws.Cells(1048576, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
  $typeVal = $ws.Cells.Item($r, 4).Value2
  if ($typeVal -eq "DATETIME") {
    $ws.Cells.Item($r, 3).Value = $typeVal
    $ws.Cells.Item($r, 4).Value = ""
  }
}

# 3) The EventDateTime column of AutoDrive_Event_100_20201215_Sample (row 3)
#    is itself part of that table's (now composite) primary key, so it
#    keeps "PK" in the type column in addition to its DATETIME data type.
$ws.Cells.Item(3, 4).Value = "PK"

# 4) Insert two more rows right after it to register the additional primary
#    key columns (ActivityID, Attribute_101) for that same table.
$ws.Range("A4:A5").EntireRow.Insert()

$ws.Cells.Item(4, 1).Value = "AutoDrive_Event_100_20201215_Sample"
$ws.Cells.Item(4, 2).Value = "ActivityID"
$ws.Cells.Item(4, 4).Value = "PK"

$ws.Cells.Item(5, 1).Value = "AutoDrive_Event_100_20201215_Sample"
$ws.Cells.Item(5, 2).Value = "Attribute_101"
$ws.Cells.Item(5, 4).Value = "PK"

# 5) Match the author's final cursor position/selection.
$ws.Range("F6").Select()
